# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F (identical on both sheets, since the
# "全部类型" sheet aggregates the same events as "展览" plus an extra row).
$updates = @{
    2  = 8833
    3  = 8204
    4  = 142
    5  = 197
    6  = 38
    8  = 142
    9  = 152
    10 = 202
    12 = 742
    13 = 203
    14 = 5293
    15 = 66
    17 = 18
    20 = 139
    21 = 3
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
